$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.583.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "'2.301.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'316.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'103.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "'39.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "'8.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").Value = "'0.969"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "'15.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").Value = "'2.649.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "'2.301.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'42.464.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'7.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'73.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").Value = "'3.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "'277.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.49%  "
$ws.Range("E24").Value = "  +20.02%  "
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -1.46%  "
$ws.Range("E28").Value = "  +3.44%  "
$ws.Range("D29").Value = "'22.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'35.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'165.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  -2.50%  "
$ws.Range("E33").Value = "  +0.25%  "
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "'2.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.78%  "
$ws.Range("D37").Value = "'0.0371"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.22%  "
$ws.Range("D38").Value = "'4.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +2.68%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("E41").Value = "  +2.58%  "
$ws.Range("D42").Value = "'69.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").Value = "'95.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "'82.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.07%  "
$ws.Range("D47").Value = "'12.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "'113.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").Value = "'8.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "'1.592.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("D51").Value = "'5.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.51%  "
